$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I11").Value = 'sd'
$ws.Range("J11").Value = 'Statement-non-opinion'
$ws.Range("I15").Value = '%'
$ws.Range("J15").Value = 'Uninterpretable'
$ws.Range("I35").Value = '%'
$ws.Range("J35").Value = 'Uninterpretable'
$ws.Range("I42").Value = 'sd'
$ws.Range("J42").Value = 'Statement-non-opinion'
$ws.Range("I46").Value = 'b'
$ws.Range("J46").Value = 'Acknowledge (Backchannel)'
$ws.Range("I51").Value = 'sv'
$ws.Range("J51").Value = 'Statement-opinion'
$ws.Range("I77").Value = 'sd'
$ws.Range("J77").Value = 'Statement-non-opinion'
$ws.Range("I79").Value = 'sd'
$ws.Range("J79").Value = 'Statement-non-opinion'
$ws.Range("I86").Value = 'sd'
$ws.Range("J86").Value = 'Statement-non-opinion'
$ws.Range("I87").Value = 'sd'
$ws.Range("J87").Value = 'Statement-non-opinion'
$ws.Range("I107").Value = 'sv'
$ws.Range("J107").Value = 'Statement-opinion'
$ws.Range("I109").Value = 'aa'
$ws.Range("J109").Value = 'Agree/Accept'
$ws.Range("I128").Value = 'aa'
$ws.Range("J128").Value = 'Agree/Accept'
$ws.Range("I129").Value = 'aa'
$ws.Range("J129").Value = 'Agree/Accept'
$ws.Range("I132").Value = 'sv'
$ws.Range("J132").Value = 'Statement-opinion'
$ws.Range("I144").Value = 'sd'
$ws.Range("J144").Value = 'Statement-non-opinion'
$ws.Range("I150").Value = 'sd'
$ws.Range("J150").Value = 'Statement-non-opinion'
$ws.Range("I154").Value = 'sv'
$ws.Range("J154").Value = 'Statement-opinion'
$ws.Range("I161").Value = 'sv'
$ws.Range("J161").Value = 'Statement-opinion'
$ws.Range("I162").Value = 'sd'
$ws.Range("J162").Value = 'Statement-non-opinion'
$ws.Range("I164").Value = 'sd'
$ws.Range("J164").Value = 'Statement-non-opinion'
$ws.Range("I179").Value = 'sd'
$ws.Range("J179").Value = 'Statement-non-opinion'
$ws.Range("I181").Value = 'ba'
$ws.Range("J181").Value = 'Appreciation'
$ws.Range("I182").Value = 'aa'
$ws.Range("J182").Value = 'Agree/Accept'
$ws.Range("I183").Value = 'sd'
$ws.Range("J183").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I221").Value = 'sd'
$ws.Range("J221").Value = 'Statement-non-opinion'
$ws.Range("I222").Value = 'sd'
$ws.Range("J222").Value = 'Statement-non-opinion'
$ws.Range("I227").Value = '%'
$ws.Range("J227").Value = 'Uninterpretable'
$ws.Range("I228").Value = 'sv'
$ws.Range("J228").Value = 'Statement-opinion'
$ws.Range("I232").Value = 'sd'
$ws.Range("J232").Value = 'Statement-non-opinion'
$ws.Range("I234").Value = 'sv'
$ws.Range("J234").Value = 'Statement-opinion'
$ws.Range("I249").Value = 'sd'
$ws.Range("J249").Value = 'Statement-non-opinion'
$ws.Range("I250").Value = 'sd'
$ws.Range("J250").Value = 'Statement-non-opinion'
$ws.Range("I251").Value = 'sv'
$ws.Range("J251").Value = 'Statement-opinion'
$ws.Range("I261").Value = 'aa'
$ws.Range("J261").Value = 'Agree/Accept'
$ws.Range("I264").Value = 'sd'
$ws.Range("J264").Value = 'Statement-non-opinion'
$ws.Range("I266").Value = 'sd'
$ws.Range("J266").Value = 'Statement-non-opinion'
$ws.Range("I268").Value = 'aa'
$ws.Range("J268").Value = 'Agree/Accept'
$ws.Range("I272").Value = 'sd'
$ws.Range("J272").Value = 'Statement-non-opinion'
$ws.Range("I274").Value = 'sd'
$ws.Range("J274").Value = 'Statement-non-opinion'
$ws.Range("I276").Value = '%'
$ws.Range("J276").Value = 'Uninterpretable'
$ws.Range("I277").Value = 'aa'
$ws.Range("J277").Value = 'Agree/Accept'
$ws.Range("I292").Value = 'sd'
$ws.Range("J292").Value = 'Statement-non-opinion'
$ws.Range("I306").Value = '%'
$ws.Range("J306").Value = 'Uninterpretable'
$ws.Range("I307").Value = 'sd'
$ws.Range("J307").Value = 'Statement-non-opinion'
$ws.Range("I311").Value = 'sd'
$ws.Range("J311").Value = 'Statement-non-opinion'
$ws.Range("I316").Value = 'sd'
$ws.Range("J316").Value = 'Statement-non-opinion'
$ws.Range("I318").Value = 'ba'
$ws.Range("J318").Value = 'Appreciation'
$ws.Range("I324").Value = 'aa'
$ws.Range("J324").Value = 'Agree/Accept'
$ws.Range("I328").Value = 'aa'
$ws.Range("J328").Value = 'Agree/Accept'
$ws.Range("I331").Value = 'aa'
$ws.Range("J331").Value = 'Agree/Accept'
$ws.Range("I332").Value = 'sd'
$ws.Range("J332").Value = 'Statement-non-opinion'
$ws.Range("I340").Value = 'ba'
$ws.Range("J340").Value = 'Appreciation'
$ws.Range("I345").Value = '%'
$ws.Range("J345").Value = 'Uninterpretable'
$ws.Range("I361").Value = 'sd'
$ws.Range("J361").Value = 'Statement-non-opinion'
$ws.Range("I382").Value = 'sv'
$ws.Range("J382").Value = 'Statement-opinion'
$ws.Range("I387").Value = 'sd'
$ws.Range("J387").Value = 'Statement-non-opinion'
$ws.Range("I389").Value = 'aa'
$ws.Range("J389").Value = 'Agree/Accept'
$ws.Range("I391").Value = 'sd'
$ws.Range("J391").Value = 'Statement-non-opinion'
$ws.Range("I399").Value = 'aa'
$ws.Range("J399").Value = 'Agree/Accept'
$ws.Range("I403").Value = 'aa'
$ws.Range("J403").Value = 'Agree/Accept'
$ws.Range("I409").Value = 'sd'
$ws.Range("J409").Value = 'Statement-non-opinion'
$ws.Range("I412").Value = 'sd'
$ws.Range("J412").Value = 'Statement-non-opinion'
